$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column A (rows 2-85) as Text so the numeric-looking municipality
# codes are stored as strings, matching the CVE_MUN column semantics.
$ws.Range("A2:A85").NumberFormat = "@"

# Header: Municipio -> CVE_MUN
$ws.Range("A1").Value = "CVE_MUN"

# Replace each municipality name with its INEGI municipality code (CVE_MUN)
$codes = @("13001", "13002", "13003", "13004", "13005", "13006", "13007", "13008", "13010", "13011", "13013", "13012", "13014", "13015", "13017", "13018", "13019", "13016", "13009", "13020", "13021", "13022", "13023", "13024", "13025", "13026", "13027", "13028", "13029", "13030", "13031", "13032", "13033", "13040", "13034", "13035", "13037", "13051", "13038", "13039", "13041", "13042", "13043", "13044", "13045", "13048", "13047", "13049", "13050", "13036", "13052", "13053", "13046", "13054", "13055", "13056", "13057", "13058", "13059", "13060", "13061", "13062", "13063", "13064", "13065", "13067", "13068", "13069", "13070", "13071", "13072", "13073", "13074", "13075", "13076", "13077", "13066", "13078", "13079", "13080", "13081", "13082", "13083", "13084")
for ($i = 0; $i -lt $codes.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $codes[$i]
}

Write-Output "done"
